$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.83
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 2.42
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 4.15
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 11.4
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.33
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 1.8
$ws.Range("V2").Value = 1.88
$ws.Range("X2").Value = 7.5
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 12.5
$ws.Range("AA2").Value = 12
$ws.Range("AB2").Value = 20
$ws.Range("AC2").Value = 10.25
$ws.Range("AD2").Value = 5.8
$ws.Range("AH2").Value = 9.25
$ws.Range("AI2").Value = 16.5
$ws.Range("AJ2").Value = 10.5
$ws.Range("AK2").Value = 40
$ws.Range("AL2").Value = 26
$ws.Range("AN2").Value = 3.7
$ws.Range("AO2").Value = 9.25
$ws.Range("AP2").Value = 18
$ws.Range("AQ2").Value = 32
$ws.Range("AS2").Value = 250
$ws.Range("AU2").Value = 7.2
$ws.Range("AV2").Value = 65
$ws.Range("AW2").Value = 5.6
$ws.Range("AX2").Value = 21
$ws.Range("AY2").Value = 27
$ws.Range("AZ2").Value = 110
$ws.Range("BB2").Value = 350

# Row 3 updates
$ws.Range("G3").Value = 4.25
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 1.65
$ws.Range("J3").Value = 4.55
$ws.Range("K3").Value = 2.18
$ws.Range("L3").Value = 2.22
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 12.8
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.61
$ws.Range("Q3").Value = 1.75
$ws.Range("R3").Value = 1.85
$ws.Range("U3").Value = 1.78
$ws.Range("V3").Value = 1.89
$ws.Range("W3").Value = 10.25
$ws.Range("X3").Value = 19.5
$ws.Range("AA3").Value = 32
$ws.Range("AC3").Value = 10.75
$ws.Range("AG3").Value = 350
$ws.Range("AH3").Value = 6.1
$ws.Range("AI3").Value = 6.8
$ws.Range("AL3").Value = 10.75
$ws.Range("AM3").Value = 20
$ws.Range("AN3").Value = 6
$ws.Range("AP3").Value = 30
$ws.Range("AS3").Value = 400
$ws.Range("AT3").Value = 2.55
$ws.Range("AU3").Value = 7.6
$ws.Range("AV3").Value = 70
$ws.Range("AY3").Value = 18
$ws.Range("AZ3").Value = 27
